# Natmi following Dr Hou advice:
# recompute the Egf-Erbb2 LR-pair table for 3 sending/target clusters
# (ECs, FAPs, sCs) instead of 2 (FAPs, sCs), updating every stats column
# (E:T) with the new values and extending the sheet from 4 data rows
# (2x2 cluster pairs) to 9 data rows (3x3 cluster pairs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Egf"
$ws.Cells.Item(2,3).Value = "Erbb2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.050561
$ws.Cells.Item(2,8).Value = 0.151683
$ws.Cells.Item(2,9).Value = 0.1845256053410153
$ws.Cells.Item(2,10).Value = 0.1845256053410153
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 1.534538333333333
$ws.Cells.Item(2,14).Value = 4.603615
$ws.Cells.Item(2,15).Value = 0.1494637976135089
$ws.Cells.Item(2,16).Value = 0.1494637976135089
$ws.Cells.Item(2,17).Value = 0.07758779267166666
$ws.Cells.Item(2,18).Value = 0.698290134045
$ws.Cells.Item(2,19).Value = 0.02757989773119974
$ws.Cells.Item(2,20).Value = 0.02757989773119974

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Egf"
$ws.Cells.Item(3,3).Value = "Erbb2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.050561
$ws.Cells.Item(3,8).Value = 0.151683
$ws.Cells.Item(3,9).Value = 0.1845256053410153
$ws.Cells.Item(3,10).Value = 0.1845256053410153
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 4.973328333333334
$ws.Cells.Item(3,14).Value = 14.919985
$ws.Cells.Item(3,15).Value = 0.4844014146353658
$ws.Cells.Item(3,16).Value = 0.4844014146353658
$ws.Cells.Item(3,17).Value = 0.2514564538616667
$ws.Cells.Item(3,18).Value = 2.263108084755
$ws.Cells.Item(3,19).Value = 0.08938446426363503
$ws.Cells.Item(3,20).Value = 0.08938446426363501

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Egf"
$ws.Cells.Item(4,3).Value = "Erbb2"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 1
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.050561
$ws.Cells.Item(4,8).Value = 0.151683
$ws.Cells.Item(4,9).Value = 0.1845256053410153
$ws.Cells.Item(4,10).Value = 0.1845256053410153
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.75909
$ws.Cells.Item(4,14).Value = 11.27727
$ws.Cells.Item(4,15).Value = 0.3661347877511252
$ws.Cells.Item(4,16).Value = 0.3661347877511252
$ws.Cells.Item(4,17).Value = 0.19006334949
$ws.Cells.Item(4,18).Value = 1.71057014541
$ws.Cells.Item(4,19).Value = 0.06756124334618052
$ws.Cells.Item(4,20).Value = 0.06756124334618052

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Egf"
$ws.Cells.Item(5,3).Value = "Erbb2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.1101723333333333
$ws.Cells.Item(5,8).Value = 0.3305169999999999
$ws.Cells.Item(5,9).Value = 0.4020809813920896
$ws.Cells.Item(5,10).Value = 0.4020809813920896
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 1.534538333333333
$ws.Cells.Item(5,14).Value = 4.603615
$ws.Cells.Item(5,15).Value = 0.1494637976135089
$ws.Cells.Item(5,16).Value = 0.1494637976135089
$ws.Cells.Item(5,17).Value = 0.1690636687727777
$ws.Cells.Item(5,18).Value = 1.521573018955
$ws.Cells.Item(5,19).Value = 0.06009655042702834
$ws.Cells.Item(5,20).Value = 0.06009655042702834

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Egf"
$ws.Cells.Item(6,3).Value = "Erbb2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.1101723333333333
$ws.Cells.Item(6,8).Value = 0.3305169999999999
$ws.Cells.Item(6,9).Value = 0.4020809813920896
$ws.Cells.Item(6,10).Value = 0.4020809813920896
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 4.973328333333334
$ws.Cells.Item(6,14).Value = 14.919985
$ws.Cells.Item(6,15).Value = 0.4844014146353658
$ws.Cells.Item(6,16).Value = 0.4844014146353658
$ws.Cells.Item(6,17).Value = 0.5479231869161111
$ws.Cells.Item(6,18).Value = 4.931308682245
$ws.Cells.Item(6,19).Value = 0.1947685961843044
$ws.Cells.Item(6,20).Value = 0.1947685961843044

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Egf"
$ws.Cells.Item(7,3).Value = "Erbb2"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.1101723333333333
$ws.Cells.Item(7,8).Value = 0.3305169999999999
$ws.Cells.Item(7,9).Value = 0.4020809813920896
$ws.Cells.Item(7,10).Value = 0.4020809813920896
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.75909
$ws.Cells.Item(7,14).Value = 11.27727
$ws.Cells.Item(7,15).Value = 0.3661347877511252
$ws.Cells.Item(7,16).Value = 0.3661347877511252
$ws.Cells.Item(7,17).Value = 0.4141477165099999
$ws.Cells.Item(7,18).Value = 3.727329448589999
$ws.Cells.Item(7,19).Value = 0.1472158347807569
$ws.Cells.Item(7,20).Value = 0.1472158347807569

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Egf"
$ws.Cells.Item(8,3).Value = "Erbb2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.113272
$ws.Cells.Item(8,8).Value = 0.339816
$ws.Cells.Item(8,9).Value = 0.413393413266895
$ws.Cells.Item(8,10).Value = 0.413393413266895
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 1.534538333333333
$ws.Cells.Item(8,14).Value = 4.603615
$ws.Cells.Item(8,15).Value = 0.1494637976135089
$ws.Cells.Item(8,16).Value = 0.1494637976135089
$ws.Cells.Item(8,17).Value = 0.1738202260933333
$ws.Cells.Item(8,18).Value = 1.56438203484
$ws.Cells.Item(8,19).Value = 0.06178734945528087
$ws.Cells.Item(8,20).Value = 0.06178734945528087

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Egf"
$ws.Cells.Item(9,3).Value = "Erbb2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.113272
$ws.Cells.Item(9,8).Value = 0.339816
$ws.Cells.Item(9,9).Value = 0.413393413266895
$ws.Cells.Item(9,10).Value = 0.413393413266895
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 4.973328333333334
$ws.Cells.Item(9,14).Value = 14.919985
$ws.Cells.Item(9,15).Value = 0.4844014146353658
$ws.Cells.Item(9,16).Value = 0.4844014146353658
$ws.Cells.Item(9,17).Value = 0.5633388469733334
$ws.Cells.Item(9,18).Value = 5.07004962276
$ws.Cells.Item(9,19).Value = 0.2002483541874264
$ws.Cells.Item(9,20).Value = 0.2002483541874263

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Egf"
$ws.Cells.Item(10,3).Value = "Erbb2"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.113272
$ws.Cells.Item(10,8).Value = 0.339816
$ws.Cells.Item(10,9).Value = 0.413393413266895
$ws.Cells.Item(10,10).Value = 0.413393413266895
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.75909
$ws.Cells.Item(10,14).Value = 11.27727
$ws.Cells.Item(10,15).Value = 0.3661347877511252
$ws.Cells.Item(10,16).Value = 0.3661347877511252
$ws.Cells.Item(10,17).Value = 0.42579964248
$ws.Cells.Item(10,18).Value = 3.83219678232
$ws.Cells.Item(10,19).Value = 0.1513577096241878
$ws.Cells.Item(10,20).Value = 0.1513577096241878
